$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 0.6911764705882353
$ws.Range("I3").Value = 0.06947424892703863
$ws.Range("K3").Value = 75.25

$ws.Range("Q3").Value = 7
$ws.Range("R3").Value = 17
$ws.Range("S3").Value = 39
$ws.Range("T3").Value = 61
$ws.Range("U3").Value = 204
$ws.Range("V3").Value = 925
$ws.Range("W3").Value = 915
$ws.Range("X3").Value = 893
$ws.Range("Y3").Value = 871
$ws.Range("Z3").Value = 728

$ws.Range("AF3").Value = 0.992489
$ws.Range("AG3").Value = 0.98176
$ws.Range("AH3").Value = 0.958155
$ws.Range("AI3").Value = 0.934549
$ws.Range("AJ3").Value = 0.781116
